$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Ação" (action) for the first risk row: "Aceitar " -> "Mitigar"
$ws.Range("F2").Value = "Mitigar"

# Update the "Ação" (action) for the "Sem Acesso a Internet" risk row:
# "Transferir" -> "Mitigar", with a new "Como?" description
$ws.Range("F5").Value = "Mitigar"
$ws.Range("G5").Value = "Planejar hora extra durante a semana ou no fim de semana"

# Update the current selection to match the saved workbook state
$ws.Range("G2").Select()
